$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Arkusz1 -> Języki)
$ws.Name = "Języki"

# Fill in the language reference table (rows 2-9), so these shared
# strings are registered before the header strings below.
$data = @(
    @(1, "Wspólny"),
    @(2, "Mroczna mowa"),
    @(3, "Krasnoludzki"),
    @(4, "Elficki"),
    @(5, "Wysoki archaik"),
    @(6, "Trolli"),
    @(7, "Sekretne języki"),
    @(8, "Martwe języki")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

# Header row
$ws.Range("A1").Value = "value"
$ws.Range("B1").Value = "result"
$ws.Range("A1:B1").Font.Bold = $true

# Match column widths used by the picker table
$ws.Columns("A:B").ColumnWidth = 13.71

# Page setup matching the saved printer-ready layout
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
